$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Developer name
$ws.Range("C3").Value = "Ashmandeep Kaur"

# Column E (Preconditions / Method Inputs) filled top to bottom first
$ws.Range("E7").Value = "None"
$ws.Range("E8").Value = "None"
$ws.Range("E9").Value = "None"
$ws.Range("E10").Value = "None"
$ws.Range("E11").Value = "None"
$ws.Range("E12").Value = "Client instance exists"
$ws.Range("E13").Value = "Client instance exists"
$ws.Range("E14").Value = "Client instance exists"
$ws.Range("E15").Value = "Client instance exists"
$ws.Range("E16").Value = "Client instance exists"

# Column F (Method Inputs) filled top to bottom next
$ws.Range("F7").Value = "client_number=12345, first_name='Ashmandeep', last_name='Kaur', email='ashmandeepkaur@gmail.com'"
$ws.Range("F8").Value = "client_number='abc'"
$ws.Range("F9").Value = "first_name=''"
$ws.Range("F10").Value = "last_name=''"
$ws.Range("F11").Value = "email='invalidemail'"
$ws.Range("F12").Value = "None"
$ws.Range("F13").Value = "None"
$ws.Range("F14").Value = "None"
$ws.Range("F15").Value = "None"
$ws.Range("F16").Value = "None"

# Column G (Expected Result) filled top to bottom last
$ws.Range("G7").Value = "Client instance is created successfully with correct values."
$ws.Range("G8").Value = "Raises ValueError with message about invalid client number."
$ws.Range("G9").Value = "Raises ValueError with message about blank first name."
$ws.Range("G10").Value = "Raises ValueError with message about blank last name."
$ws.Range("G11").Value = "email_address is set to default value (e.g., None or a placeholder)."
$ws.Range("G12").Value = "Returns the client_number of the instance."
$ws.Range("G13").Value = "Returns the first_name of the instance."
$ws.Range("G14").Value = "Returns the last_name of the instance."
$ws.Range("G15").Value = "Returns the email_address of the instance."
$ws.Range("G16").Value = "Returns client details in formatted string (e.g., Client: Ashmandeep Kaur (12345))."

# Select G16 to mirror final active cell selection
$ws.Range("G16").Select()
